# Generate Report for Handback
#
# f36d11f1-0539-48c7-b681-e67590560f73.md has now been handed back (its
# localization round-trip finished), so the generated status report is
# refreshed: that file's row moves to the front of the "in-flight" block
# (right after the already-handed-back files) with its status switched to
# "Handed back: in sync with en-US" and its handback file/time filled in,
# while 11dcff0e-746e-4bde-b7da-d7fb64b12d1d.md and
# b3da3f9c-c1fd-41ef-8dc5-0c5ee32270ec.md shift down a row each, keeping
# their existing data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: rows 5-7 re-ordered (f36d11f1, then 11dcff0e, then
# b3da3f9c) with the Status columns (zh-cn / de-de) updated for the rows
# that moved.
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A5").Value = "f36d11f1-0539-48c7-b681-e67590560f73.md"
$ov.Range("B5").Value = "Handed back: in sync with en-US"
$ov.Range("C5").Value = "Handed back: in sync with en-US"

$ov.Range("A6").Value = "11dcff0e-746e-4bde-b7da-d7fb64b12d1d.md"
$ov.Range("B6").Value = "In Translation"
$ov.Range("C6").Value = "In Translation"

$ov.Range("A7").Value = "b3da3f9c-c1fd-41ef-8dc5-0c5ee32270ec.md"
$ov.Range("B7").Value = "Ready for handoff"
$ov.Range("C7").Value = "Ready for handoff"

# ---------------------------------------------------------------------
# Per-language detail sheets (zh-cn / de-de): same row re-order, plus the
# freshly-handed-back file's Latest Target File / Latest Handback File /
# Latest Handback DateTime columns (E, F, G) get populated, and the rows
# that shift down (11dcff0e, b3da3f9c) keep E/F blank as before.
# ---------------------------------------------------------------------
$langs = @(
    @{ Name = "zh-cn";
       F36Xlf = "f36d11f1-0539-48c7-b681-e67590560f73.1333ae907a1131cc8e7aa38273f30094cd4dd266.zh-cn.xlf";
       F36D   = "2016-03-10 12:30:30";
       F36G   = "2016-03-10 12:31:14";
       Dc0Xlf = "11dcff0e-746e-4bde-b7da-d7fb64b12d1d.34e508a11b7473be85ca32bf46b6b9b53099aebc.zh-cn.xlf";
       Dc0D   = "2016-03-10 12:26:41";
       B3dXlf = "b3da3f9c-c1fd-41ef-8dc5-0c5ee32270ec.3318fbd1b21412a676b090a4f77a66da62003549.zh-cn.xlf";
       B3dD   = "2016-03-10 12:30:30";
     },
    @{ Name = "de-de";
       F36Xlf = "f36d11f1-0539-48c7-b681-e67590560f73.1333ae907a1131cc8e7aa38273f30094cd4dd266.de-de.xlf";
       F36D   = "2016-03-10 12:30:40";
       F36G   = "2016-03-10 12:31:27";
       Dc0Xlf = "11dcff0e-746e-4bde-b7da-d7fb64b12d1d.34e508a11b7473be85ca32bf46b6b9b53099aebc.de-de.xlf";
       Dc0D   = "2016-03-10 12:27:17";
       B3dXlf = "b3da3f9c-c1fd-41ef-8dc5-0c5ee32270ec.3318fbd1b21412a676b090a4f77a66da62003549.de-de.xlf";
       B3dD   = "2016-03-10 12:30:40";
     }
)

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.Name)

    # Row 5: f36d11f1 - now fully handed back.
    $ws.Range("A5").Value = "f36d11f1-0539-48c7-b681-e67590560f73.md"
    $ws.Range("B5").Value = "Handed back: in sync with en-US"
    $ws.Range("C5").Value = $lang.F36Xlf
    $ws.Range("D5").Value = $lang.F36D
    $ws.Range("E5").Value = "f36d11f1-0539-48c7-b681-e67590560f73.md"
    $ws.Range("F5").Value = $lang.F36Xlf
    $ws.Range("G5").Value = $lang.F36G
    $ws.Range("H5").Value = "Include"

    # Row 6: 11dcff0e - still in translation, shifted down from row 5.
    $ws.Range("A6").Value = "11dcff0e-746e-4bde-b7da-d7fb64b12d1d.md"
    $ws.Range("B6").Value = "In Translation"
    $ws.Range("C6").Value = $lang.Dc0Xlf
    $ws.Range("D6").Value = $lang.Dc0D
    $ws.Range("E6").ClearContents()
    $ws.Range("F6").ClearContents()
    $ws.Range("G6").Value = "0001-01-01 00:00:00"
    $ws.Range("H6").Value = "Include"

    # Row 7: b3da3f9c - still ready for handoff, shifted down from row 6.
    $ws.Range("A7").Value = "b3da3f9c-c1fd-41ef-8dc5-0c5ee32270ec.md"
    $ws.Range("B7").Value = "Ready for handoff"
    $ws.Range("C7").Value = $lang.B3dXlf
    $ws.Range("D7").Value = $lang.B3dD
    $ws.Range("E7").ClearContents()
    $ws.Range("F7").ClearContents()
    $ws.Range("G7").Value = "0001-01-01 00:00:00"
    $ws.Range("H7").Value = "Include"
}
